# fix: revert admin dev default; seed customers only when table empty;
# autosave on customer select when hours/day present
#
# This reconstructs the weekly export for Chris Zavesky for 2026-01-19:
#  - day 1 becomes a PTO entry (6.5h, rate/total reset to 0 - admin default reverted)
#  - days 2-5 become fresh "seeded" customers (Richer/Durfee/Tercek/Patton) on
#    2026-01-20..2026-01-23, all Regular, rate/total reset to 0
#  - the hourly SUBTOTAL row total drops from 42 -> 32 and its Reg:/OT: label
#    + dollar total follow along (reverted to the non-admin default of 0)
#  - three new rows are appended below the ADMIN separator: HOURLY SUBTOTAL,
#    ADMIN SUBTOTAL and GRAND TOTAL
#  - the "Jason Schema" autosave mirror sheet is regenerated the same way,
#    plus the employee id is rotated to the new short-form id

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Weekly Timesheet"
$ws2 = $wb.Worksheets.Item(2)   # "Jason Schema"

# ---------------------------------------------------------------------------
# 1) Weekly Timesheet - rows 2-6 (row-major, left-to-right, so shared strings
#    land in the same first-use order the exporter would produce them in)
# ---------------------------------------------------------------------------

# Row 2: 2026-01-19 PTO day - date (A2) is unchanged
$ws1.Cells.Item(2,2).Value = "PTO"
$ws1.Cells.Item(2,3).Value = 6.5
$ws1.Cells.Item(2,4).Value = "PTO"
$ws1.Cells.Item(2,5).Value = 0
$ws1.Cells.Item(2,6).Value = 0

# Row 3: 2026-01-20 / Richer / 7h / Regular
$ws1.Cells.Item(3,2).Value = "Richer"
$ws1.Cells.Item(3,5).Value = 0
$ws1.Cells.Item(3,6).Value = 0

# Row 4: 2026-01-21 / Durfee / 6.5h / Regular
$ws1.Cells.Item(4,2).Value = "Durfee"
$ws1.Cells.Item(4,3).Value = 6.5
$ws1.Cells.Item(4,5).Value = 0
$ws1.Cells.Item(4,6).Value = 0

# Row 5: 2026-01-22 / Tercek / 6h / Regular
$ws1.Cells.Item(5,2).Value = "Tercek"
$ws1.Cells.Item(5,3).Value = 6
$ws1.Cells.Item(5,5).Value = 0
$ws1.Cells.Item(5,6).Value = 0

# Row 6: date moves 2026-01-25 -> 2026-01-23 / Patton / 6h / Regular
$ws1.Cells.Item(6,1).Value = "2026-01-23"
$ws1.Cells.Item(6,2).Value = "Patton"
$ws1.Cells.Item(6,3).Value = 6
$ws1.Cells.Item(6,5).Value = 0
$ws1.Cells.Item(6,6).Value = 0

# ---------------------------------------------------------------------------
# 2) Weekly Timesheet - SUBTOTAL row (8): hours 42 -> 32, label + $ total
#    revert to the non-admin default (0) - A8/B8/E8 (labels / blanks) are
#    untouched so their existing shared-string cells survive as-is.
# ---------------------------------------------------------------------------
$ws1.Cells.Item(8,3).Value = 32
$ws1.Cells.Item(8,4).Value = "Reg: 32 / OT: 0"
$ws1.Cells.Item(8,6).Value = 0

# ---------------------------------------------------------------------------
# 3) Weekly Timesheet - three brand-new summary rows under the ADMIN
#    separator (row 9 is left completely untouched)
# ---------------------------------------------------------------------------
$ws1.Cells.Item(11,1).Value = "HOURLY SUBTOTAL"
$ws1.Cells.Item(11,2).Value = ""
$ws1.Cells.Item(11,3).Value = ""
$ws1.Cells.Item(11,4).Value = ""
$ws1.Cells.Item(11,5).Value = ""
$ws1.Cells.Item(11,6).Value = 0

$ws1.Cells.Item(12,1).Value = "ADMIN SUBTOTAL"
$ws1.Cells.Item(12,2).Value = ""
$ws1.Cells.Item(12,3).Value = ""
$ws1.Cells.Item(12,4).Value = ""
$ws1.Cells.Item(12,5).Value = ""
$ws1.Cells.Item(12,6).Value = 0

$ws1.Cells.Item(13,1).Value = "GRAND TOTAL"
$ws1.Cells.Item(13,2).Value = ""
$ws1.Cells.Item(13,3).Value = ""
$ws1.Cells.Item(13,4).Value = ""
$ws1.Cells.Item(13,5).Value = ""
$ws1.Cells.Item(13,6).Value = 0

# Style row 11 + 12 (bold font on the light-tan "hourly/admin subtotal" fill)
# by stamping the exact format of an already-bold/filled cell (format-painter
# style) so the underlying font/fill objects get reused intact.
$ws1.Cells.Item(8,1).Copy()
$ws1.Range("A11:D12").PasteSpecial(-4122)
$ws1.Cells.Item(8,5).Copy()
$ws1.Range("E11:F12").PasteSpecial(-4122)

# Style row 13 (bold red font on the light-green "grand total" fill) - built
# fresh once on A13 then stamped across the rest of the row.
$ws1.Cells.Item(13,1).Font.Bold = $true
$ws1.Cells.Item(13,1).Font.Color = 255
$ws1.Cells.Item(13,1).Interior.Color = 0xE0F8E8
$ws1.Cells.Item(13,1).Copy()
$ws1.Range("B13:D13").PasteSpecial(-4122)

$ws1.Cells.Item(13,6).NumberFormat = '"$"#,##0.00'
$ws1.Cells.Item(13,6).Font.Bold = $true
$ws1.Cells.Item(13,6).Font.Color = 255
$ws1.Cells.Item(13,6).Interior.Color = 0xE0F8E8
$ws1.Cells.Item(13,6).Copy()
$ws1.Cells.Item(13,5).PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Jason Schema (autosave mirror) - regenerate rows 2-6 the same way, plus
#    rotate the employee id to the new short id
# ---------------------------------------------------------------------------

# Row 2
$ws2.Cells.Item(2,2).Value = "emp_5chpvt65"
$ws2.Cells.Item(2,4).Value = "PTO"
$ws2.Cells.Item(2,5).Value = 6.5
$ws2.Cells.Item(2,6).Value = 0
$ws2.Cells.Item(2,7).Value = 0
$ws2.Cells.Item(2,8).Value = "PTO"
$ws2.Cells.Item(2,9).Value = "PTO"

# Row 3
$ws2.Cells.Item(3,2).Value = "emp_5chpvt65"
$ws2.Cells.Item(3,4).Value = "Richer"
$ws2.Cells.Item(3,6).Value = 0
$ws2.Cells.Item(3,7).Value = 0

# Row 4
$ws2.Cells.Item(4,2).Value = "emp_5chpvt65"
$ws2.Cells.Item(4,4).Value = "Durfee"
$ws2.Cells.Item(4,5).Value = 6.5
$ws2.Cells.Item(4,6).Value = 0
$ws2.Cells.Item(4,7).Value = 0

# Row 5
$ws2.Cells.Item(5,2).Value = "emp_5chpvt65"
$ws2.Cells.Item(5,4).Value = "Tercek"
$ws2.Cells.Item(5,5).Value = 6
$ws2.Cells.Item(5,6).Value = 0
$ws2.Cells.Item(5,7).Value = 0

# Row 6
$ws2.Cells.Item(6,2).Value = "emp_5chpvt65"
$ws2.Cells.Item(6,3).Value = "2026-01-23"
$ws2.Cells.Item(6,4).Value = "Patton"
$ws2.Cells.Item(6,5).Value = 6
$ws2.Cells.Item(6,6).Value = 0
$ws2.Cells.Item(6,7).Value = 0

$wb.Save()
